# ---------------------------------------------------------------------
# CompStat weekly refresh: "New crime data collected"
#
# Rolls the report forward one week:
#   Volume 31 Number 2            -> Volume 31 Number 3
#   Week of 1/8/2024 - 1/14/2024   -> Week of 1/15/2024 - 1/21/2024
# and refreshes every statistic in the main precinct comparison table
# (rows 14-30, columns C:N) with the newly collected counts / % changes.
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume number & reporting week -----------------------------
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Rows 28-30 (Shooting Vic. / Shooting Inc. / Hate Crimes): a few cells
# move from the "no data yet" text placeholders ("0" / "***.*") to real
# numeric figures now that data exists. Give them the same number formats
# their numeric neighbours already use before writing the values, so they
# land on the existing #,##0 / #,##0.0 styles instead of minting new ones.
$fmtCount = "#,##0"
$fmtPct = "#,##0.0;""-""#,##0.0"

$ws.Range("C28").NumberFormat = $fmtCount
$ws.Range("F28").NumberFormat = $fmtCount
$ws.Range("I28").NumberFormat = $fmtCount
$ws.Range("C29").NumberFormat = $fmtCount
$ws.Range("F29").NumberFormat = $fmtCount
$ws.Range("I29").NumberFormat = $fmtCount
$ws.Range("D30").NumberFormat = $fmtCount
$ws.Range("G30").NumberFormat = $fmtCount
$ws.Range("J30").NumberFormat = $fmtCount

$ws.Range("E30").NumberFormat = $fmtPct
$ws.Range("H30").NumberFormat = $fmtPct
$ws.Range("K30").NumberFormat = $fmtPct

# --- Newly collected figures (rows 14-30, columns C:N) -----------------
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = -50
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = -50
$ws.Range("N14").Value = -88.888888888888
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 15
$ws.Range("H15").Value = -20
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -61.538461538461
$ws.Range("L15").Value = -54.545454545454
$ws.Range("M15").Value = -37.5
$ws.Range("N15").Value = -68.75
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 32
$ws.Range("E16").Value = 3.125
$ws.Range("F16").Value = 170
$ws.Range("G16").Value = 121
$ws.Range("H16").Value = 40.495867768595
$ws.Range("I16").Value = 114
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = 18.75
$ws.Range("L16").Value = 48.051948051948
$ws.Range("M16").Value = -15.555555555555
$ws.Range("N16").Value = -82.352941176470
$ws.Range("C17").Value = 35
$ws.Range("D17").Value = 53
$ws.Range("E17").Value = -33.962264150943
$ws.Range("F17").Value = 211
$ws.Range("G17").Value = 192
$ws.Range("H17").Value = 9.895833333333
$ws.Range("I17").Value = 148
$ws.Range("J17").Value = 148
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 15.625
$ws.Range("M17").Value = 89.743589743589
$ws.Range("N17").Value = -8.074534161490
$ws.Range("C18").Value = 43
$ws.Range("D18").Value = 50
$ws.Range("E18").Value = -14
$ws.Range("F18").Value = 148
$ws.Range("H18").Value = -12.426035502958
$ws.Range("I18").Value = 112
$ws.Range("J18").Value = 127
$ws.Range("K18").Value = -11.811023622047
$ws.Range("L18").Value = 14.285714285714
$ws.Range("M18").Value = -34.117647058823
$ws.Range("N18").Value = -89.230769230769
$ws.Range("C19").Value = 123
$ws.Range("D19").Value = 110
$ws.Range("E19").Value = 11.818181818181
$ws.Range("F19").Value = 468
$ws.Range("G19").Value = 458
$ws.Range("H19").Value = 2.183406113537
$ws.Range("I19").Value = 355
$ws.Range("J19").Value = 348
$ws.Range("K19").Value = 2.011494252873
$ws.Range("L19").Value = -31.992337164751
$ws.Range("M19").Value = 38.671875
$ws.Range("N19").Value = -19.501133786848
$ws.Range("C20").Value = 53
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = 60.606060606060
$ws.Range("F20").Value = 182
$ws.Range("G20").Value = 165
$ws.Range("H20").Value = 10.303030303030
$ws.Range("I20").Value = 135
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = 13.445378151260
$ws.Range("L20").Value = 35
$ws.Range("M20").Value = 46.739130434782
$ws.Range("N20").Value = -91.041804910418
$ws.Range("C21").Value = 288
$ws.Range("D21").Value = 282
$ws.Range("E21").Value = 2.127659574468
$ws.Range("F21").Value = 1192
$ws.Range("G21").Value = 1122
$ws.Range("H21").Value = 6.238859180035
$ws.Range("I21").Value = 870
$ws.Range("J21").Value = 853
$ws.Range("K21").Value = 1.992966002344
$ws.Range("L21").Value = -7.348242811501
$ws.Range("M21").Value = 17.567567567567
$ws.Range("N21").Value = -77.225130890052
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = 17.647058823529
$ws.Range("I22").Value = 14
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 55.555555555555
$ws.Range("L22").Value = -30
$ws.Range("M22").Value = -6.666666666666
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 22
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -12
$ws.Range("I23").Value = 13
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = -13.333333333333
$ws.Range("L23").Value = -7.142857142857
$ws.Range("M23").Value = 30
$ws.Range("C24").Value = 347
$ws.Range("D24").Value = 294
$ws.Range("E24").Value = 18.027210884353
$ws.Range("F24").Value = 1310
$ws.Range("G24").Value = 1152
$ws.Range("H24").Value = 13.715277777777
$ws.Range("I24").Value = 970
$ws.Range("J24").Value = 857
$ws.Range("K24").Value = 13.185530921820
$ws.Range("L24").Value = 26.137841352405
$ws.Range("M24").Value = 105.944798301486
$ws.Range("C25").Value = 87
$ws.Range("D25").Value = 104
$ws.Range("E25").Value = -16.346153846153
$ws.Range("F25").Value = 380
$ws.Range("G25").Value = 370
$ws.Range("H25").Value = 2.702702702702
$ws.Range("I25").Value = 265
$ws.Range("J25").Value = 280
$ws.Range("K25").Value = -5.357142857142
$ws.Range("L25").Value = 11.814345991561
$ws.Range("M25").Value = 6.425702811244
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 22
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -27.777777777777
$ws.Range("L26").Value = -18.75
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = -60
$ws.Range("F27").Value = 32
$ws.Range("G27").Value = 36
$ws.Range("H27").Value = -11.111111111111
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = -10.344827586206
$ws.Range("L27").Value = 13.043478260869
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -85.714285714285
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = -83.333333333333
$ws.Range("L28").Value = -75
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -94.444444444444
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -66.666666666666
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = -83.333333333333
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = -80
$ws.Range("L29").Value = -75
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -94.444444444444
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 100
